$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously-used range (A1:J4) so stale columns G:J / rows 3:4 disappear.
$ws.Range("A1:J4").Clear()

# New header row (1) -- six columns of labels.
$ws.Range("A1").Value = "strain_crit_rup_casing"
$ws.Range("B1").Value = "strain_crit_rup_tubing"
$ws.Range("C1").Value = "sigma_strain_crit_rup_casing"
$ws.Range("D1").Value = "sigma_strain_crit_rup_tubing"
$ws.Range("E1").Value = "sigma_mu_strain_crit_rup_casing"
$ws.Range("F1").Value = "sigma_mu_strain_crit_rup_tubing"

# New data row (2) -- literal (already-computed) numeric values.
$ws.Range("A2").Value = 58.75371104
$ws.Range("B2").Value = 54.89140057
$ws.Range("C2").Value = 0.185657455
$ws.Range("D2").Value = 0.392014819
$ws.Range("E2").Value = 0.102984227
$ws.Range("F2").Value = 0.261343213

# Reset column widths back to default (remove the bespoke widths from before).
$ws.Columns.Item(1).ColumnWidth = 8.43
$ws.Columns.Item(2).ColumnWidth = 8.43
$ws.Columns.Item(3).ColumnWidth = 8.43
$ws.Columns.Item(4).ColumnWidth = 8.43

$ws.Range("F5").Select()
